# Add "NA" values under the duplicate_image_filename column (column E)
# for the practice rows (2-5) and trial rows (6-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
